# Update column F ("dSF") values on the active worksheet to reflect
# repulled / recalculated data, per commit: "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F5").Value = 11
$ws.Range("F6").Value = 6
$ws.Range("F10").Value = 0
$ws.Range("F14").Value = -1
$ws.Range("F20").Value = 2
